$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- Shared-string ordering matters: the order in which *new* unique string
# values are introduced determines the slot they land in when the workbook is
# saved (the orphaned slot left behind by the old "Mise en place" string, then
# new slots appended in first-seen order). This sequence reproduces:
#   75 = "Implémentation de la base mvc du site web."
#   76 = "Enivronnement"
#   77 = "Implémentation"
#   78 = "Mise en place du lien entre le site et la base de données."
$ws.Range("H56").Value = "Implémentation de la base mvc du site web."
$ws.Range("G58").Value = "Enivronnement"
$ws.Range("G56").Value = "Implémentation"
$ws.Range("H58").Value = "Mise en place du lien entre le site et la base de données."

# --- Row 57: brand-new data row ---
$ws.Range("A57").Value = 44326
$ws.Range("B57").Value = 2
$ws.Range("C57").Value = 0.63888888888888895
$ws.Range("D57").Value = 0.68958333333333333
$ws.Range("E57").Formula = "=D57-C57"
$ws.Range("F57").Value = "Réalisation"
$ws.Range("G57").Value = "Implémentation"
$ws.Range("H57").Value = "Implémentation de la base mvc du site web."

# --- Row 58: brand-new data (G58/H58 string values already set above) ---
$ws.Range("A58").Value = 44326
$ws.Range("B58").Value = 2
$ws.Range("C58").Value = 0.68958333333333333
$ws.Range("F58").Value = "Réalisation"

# --- Row 56 remaining numeric/formula updates ---
$ws.Range("D56").Value = 0.62847222222222221
$ws.Range("E56").Formula = "=D56-C56"

# --- Row height / column width / selection cosmetics ---
$ws.Rows.Item(58).RowHeight = 30
$ws.Columns.Item(8).ColumnWidth = 40.5
$ws.Range("A59").Select()
